$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text strings (e.g. thousand-separated "36.462.99"),
# some of which look numeric (e.g. "245.79"). Force text format on the cells
# whose price is being refreshed so COM does not coerce them to Number.
foreach ($addr in @("D2","D3","D5","D7","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D21","D23","D25","D26","D27","D28","D30","D32","D33","D40","D42","D45","D46","D47","D48","D49","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '36.462.99'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '1.978.73'
$ws.Range("E3").Value = '  -3.28%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '245.79'
$ws.Range("E5").Value = '  -2.77%  '
$ws.Range("E6").Value = '  -4.35%  '
$ws.Range("D7").Value = '59.10'
$ws.Range("E7").Value = '  -9.86%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.374'
$ws.Range("E9").Value = '  -8.63%  '
$ws.Range("D10").Value = '56.67'
$ws.Range("E10").Value = '  -5.13%  '
$ws.Range("D11").Value = '0.0856'
$ws.Range("E11").Value = '  +8.54%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '22.74'
$ws.Range("E13").Value = '  -3.80%  '
$ws.Range("D14").Value = '0.859'
$ws.Range("E14").Value = '  -7.15%  '
$ws.Range("D15").Value = '2.266.90'
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").Value = '13.79'
$ws.Range("E16").Value = '  -6.95%  '
$ws.Range("E17").Value = '  -5.48%  '
$ws.Range("D18").Value = '1.972.83'
$ws.Range("E18").Value = '  -3.66%  '
$ws.Range("D19").Value = '36.373.25'
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("D20").Value = '0.0₃0895'
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").Value = '70.46'
$ws.Range("E21").Value = '  -4.31%  '
$ws.Range("E22").Value = '  -5.06%  '
$ws.Range("D23").Value = '233.99'
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '2.50'
$ws.Range("E25").Value = '  -4.79%  '
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").Value = '9.87'
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("D28").Value = '163.09'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("D30").Value = '19.85'
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").Value = '1.19'
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").Value = '4.88'
$ws.Range("E33").Value = '  -5.88%  '
$ws.Range("E34").Value = '  +4.53%  '
$ws.Range("E35").Value = '  -5.39%  '
$ws.Range("E36").Value = '  -3.56%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("E38").Value = '  -1.34%  '
$ws.Range("E39").Value = '  -7.14%  '
$ws.Range("D40").Value = '2.91'
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("E41").Value = '  -4.96%  '
$ws.Range("D42").Value = '0.0964'
$ws.Range("E42").Value = '  -5.56%  '
$ws.Range("E43").Value = '  -5.68%  '
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D45").Value = '1.09'
$ws.Range("E45").Value = '  -6.74%  '
$ws.Range("D46").Value = '16.17'
$ws.Range("E46").Value = '  -9.16%  '
$ws.Range("D47").Value = '91.31'
$ws.Range("E47").Value = '  -5.40%  '
$ws.Range("D48").Value = '1.370.40'
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("D49").Value = '7.45'
$ws.Range("E49").Value = '  -5.52%  '
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("D51").Value = '45.43'
$ws.Range("E51").Value = '  -4.97%  '
